$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02032019")

$ws.Range("F4").Value = 37.875087999999998
$ws.Range("G4").Value = -122.260554

$ws.Range("F11").Value = 37.883595
$ws.Range("G11").Value = -122.303265

$ws.Range("G12").Select()
